$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4725.0117
$ws.Range("I15").Value = 4725.0117
$ws.Range("K15").Value = 14175.0351
$ws.Range("M15").Value = -14006.0351

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1180.6666
$ws.Range("J129").Value = 1033.5883
$ws.Range("L129").Value = 3100.7649
$ws.Range("N129").Value = -13100.7649

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 1765.138
$ws.Range("I131").Value = 1116.1333
$ws.Range("J131").Value = 2460.5
$ws.Range("K131").Value = 3348.3999
$ws.Range("L131").Value = 7381.5
$ws.Range("M131").Value = 1691.6001
$ws.Range("N131").Value = -17461.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4062557
$ws.Range("I137").Value = 12843713
$ws.Range("J137").Value = 9716
$ws.Range("K137").Value = 38531139
$ws.Range("L137").Value = 29148
$ws.Range("M137").Value = -38528589
$ws.Range("N137").Value = -34248

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2499.1936
$ws.Range("I138").Value = 1412.04
$ws.Range("J138").Value = 3233.7568
$ws.Range("K138").Value = 4236.12
$ws.Range("L138").Value = 9701.270400000001
$ws.Range("M138").Value = 903.8800000000001
$ws.Range("N138").Value = -19981.2704

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2077.476
$ws.Range("I2").Value = 2360.0667
$ws.Range("J2").Value = 1371
$ws.Range("K2").Value = 2360.0667
$ws.Range("L2").Value = 1371
$ws.Range("M2").Value = -2247.0667
$ws.Range("N2").Value = -1597

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1799.0278
$ws.Range("I61").Value = 1160.4615
$ws.Range("J61").Value = 3459.3
$ws.Range("K61").Value = 1160.4615
$ws.Range("L61").Value = 3459.3
$ws.Range("M61").Value = -948.4614999999999
$ws.Range("N61").Value = -3883.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 11617292
$ws.Range("I88").Value = 50004476
$ws.Range("J88").Value = 2020496.2
$ws.Range("K88").Value = 50004476
$ws.Range("L88").Value = 2020496.2
$ws.Range("M88").Value = -50004070
$ws.Range("N88").Value = -2021308.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 11617292
$ws.Range("I91").Value = 50004476
$ws.Range("J91").Value = 2020496.2
$ws.Range("K91").Value = 50004476
$ws.Range("L91").Value = 2020496.2
$ws.Range("M91").Value = -50003072
$ws.Range("N91").Value = -2023304.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2077.476
$ws.Range("I116").Value = 2360.0667
$ws.Range("J116").Value = 1371
$ws.Range("K116").Value = 2360.0667
$ws.Range("L116").Value = 1371
$ws.Range("M116").Value = -66.06669999999986
$ws.Range("N116").Value = -5959

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1799.0278
$ws.Range("I136").Value = 1160.4615
$ws.Range("J136").Value = 3459.3
$ws.Range("K136").Value = 3481.3845
$ws.Range("L136").Value = 10377.9
$ws.Range("M136").Value = -931.3844999999997
$ws.Range("N136").Value = -15477.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7266.2
$ws.Range("I20").Value = 1943.8889
$ws.Range("J20").Value = 11620.818
$ws.Range("K20").Value = 1943.8889
$ws.Range("L20").Value = 11620.818
$ws.Range("M20").Value = -1696.8889
$ws.Range("N20").Value = -12114.818

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3750.8333
$ws.Range("I86").Value = 2501.2
$ws.Range("J86").Value = 9999
$ws.Range("K86").Value = 2501.2
$ws.Range("L86").Value = 9999
$ws.Range("M86").Value = -1378.2
$ws.Range("N86").Value = -12245

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3750.8333
$ws.Range("I89").Value = 2501.2
$ws.Range("J89").Value = 9999
$ws.Range("K89").Value = 12506
$ws.Range("L89").Value = 49995
$ws.Range("M89").Value = -6890
$ws.Range("N89").Value = -61227

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3430.1738
$ws.Range("I105").Value = 3018.3635
$ws.Range("J105").Value = 3807.6667
$ws.Range("K105").Value = 3018.3635
$ws.Range("L105").Value = 3807.6667
$ws.Range("M105").Value = -1271.3635
$ws.Range("N105").Value = -7301.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3370102.5
$ws.Range("I31").Value = 1365.36
$ws.Range("J31").Value = 6807589.5
$ws.Range("K31").Value = 1365.36
$ws.Range("L31").Value = 6807589.5
$ws.Range("M31").Value = -1070.36
$ws.Range("N31").Value = -6808179.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3370102.5
$ws.Range("I34").Value = 1365.36
$ws.Range("J34").Value = 6807589.5
$ws.Range("K34").Value = 1365.36
$ws.Range("L34").Value = 6807589.5
$ws.Range("M34").Value = -1163.36
$ws.Range("N34").Value = -6807993.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2059.125
$ws.Range("I99").Value = 2166
$ws.Range("J99").Value = 1968.6923
$ws.Range("K99").Value = 2166
$ws.Range("L99").Value = 1968.6923
$ws.Range("M99").Value = -668
$ws.Range("N99").Value = -4964.6923

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H119").Value = 58756.555
$ws.Range("J119").Value = 58756.555
$ws.Range("L119").Value = 58756.555
$ws.Range("N119").Value = -68432.55499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2059.125
$ws.Range("I126").Value = 2166
$ws.Range("J126").Value = 1968.6923
$ws.Range("K126").Value = 6498
$ws.Range("L126").Value = 5906.0769
$ws.Range("M126").Value = -4028
$ws.Range("N126").Value = -10846.0769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 743.38464
$ws.Range("I23").Value = 418
$ws.Range("J23").Value = 1475.5
$ws.Range("K23").Value = 1254
$ws.Range("L23").Value = 4426.5
$ws.Range("M23").Value = -1019
$ws.Range("N23").Value = -4896.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1244.1647
$ws.Range("I68").Value = 969.9091
$ws.Range("J68").Value = 1339.9365
$ws.Range("K68").Value = 2909.7273
$ws.Range("L68").Value = 4019.8095
$ws.Range("M68").Value = -2098.7273
$ws.Range("N68").Value = -5641.8095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1244.1647
$ws.Range("I71").Value = 969.9091
$ws.Range("J71").Value = 1339.9365
$ws.Range("K71").Value = 8729.1819
$ws.Range("L71").Value = 12059.4285
$ws.Range("M71").Value = -4673.1819
$ws.Range("N71").Value = -20171.4285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5163.977
$ws.Range("I70").Value = 5139.39
$ws.Range("J70").Value = 5500
$ws.Range("K70").Value = 5139.39
$ws.Range("L70").Value = 5500
$ws.Range("M70").Value = -4869.39
$ws.Range("N70").Value = -6040

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5163.977
$ws.Range("I73").Value = 5139.39
$ws.Range("J73").Value = 5500
$ws.Range("K73").Value = 5139.39
$ws.Range("L73").Value = 5500
$ws.Range("M73").Value = -4203.39
$ws.Range("N73").Value = -7372

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3254.5186
$ws.Range("I97").Value = 2429.0908
$ws.Range("J97").Value = 6886.4
$ws.Range("K97").Value = 2429.0908
$ws.Range("L97").Value = 6886.4
$ws.Range("M97").Value = -1933.0908
$ws.Range("N97").Value = -7878.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1337.3334
$ws.Range("I102").Value = 1006
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1006
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 616
$ws.Range("N102").Value = -5244

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6067.4614
$ws.Range("I126").Value = 14438
$ws.Range("K126").Value = 43314
$ws.Range("M126").Value = -40844

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 33337490
$ws.Range("I132").Value = 52635590
$ws.Range("J132").Value = 4402.4546
$ws.Range("K132").Value = 157906770
$ws.Range("L132").Value = 13207.3638
$ws.Range("M132").Value = -157904240
$ws.Range("N132").Value = -18267.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2106.6875
$ws.Range("I61").Value = 1830.2
$ws.Range("J61").Value = 2567.5
$ws.Range("K61").Value = 1830.2
$ws.Range("L61").Value = 2567.5
$ws.Range("M61").Value = -1628.2
$ws.Range("N61").Value = -2971.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2106.6875
$ws.Range("I113").Value = 1830.2
$ws.Range("J113").Value = 2567.5
$ws.Range("K113").Value = 1830.2
$ws.Range("L113").Value = 2567.5
$ws.Range("M113").Value = 339.8
$ws.Range("N113").Value = -6907.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 55068.844
$ws.Range("I122").Value = 73679.14
$ws.Range("J122").Value = 2960
$ws.Range("K122").Value = 221037.42
$ws.Range("L122").Value = 8880
$ws.Range("M122").Value = -218587.42
$ws.Range("N122").Value = -13780

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 66003.8
$ws.Range("J31").Value = 66003.8
$ws.Range("L31").Value = 66003.8
$ws.Range("N31").Value = -66699.8
